# "added frame for leaderboard"
# Append a new leaderboard row (janedoey) to the Users sheet, right after
# the existing rows, extending the used range from A1:B3 to A1:B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "janedoey"
$ws.Range("B4").Value = "c78269a8b5134b8f79ae1f2dbb124979effcd3dcad53abf3d4e170b835823847"
